# Rename the three inline logo pictures that live in the document's
# headers/footers:
#   - footer (default)  -> Pearson logo  : image1.png -> image2.png
#   - footer (first page)-> Pearson logo  : image1.png -> image2.png
#   - header (first page)-> BTec logo     : image2.jpg -> image1.jpg
#
# Note: InlineShapes that are not the very first paragraph of their
# header/footer story can't be renamed reliably by holding on to the
# InlineShape reference returned from HeaderFooter.Range.InlineShapes
# directly (the freshly-fetched object's address goes stale as soon as
# you try to write to it). Selecting the shape's range first and then
# reaching the shape back through $word.Selection re-resolves the
# address correctly, so we always go through the Selection object.

$d = $word.ActiveDocument

function Rename-InlineLogo($headerFooter, $newName) {
    $shape = $headerFooter.Range.InlineShapes.Item(1)
    $shape.Range.Select()
    $word.Selection.InlineShapes.Item(1).Name = $newName
}

$section = $d.Sections.Item(1)

# Pearson Edexcel logo, default footer
Rename-InlineLogo $section.Footers.Item(1) "image2.png"

# Pearson Edexcel logo, first-page footer
Rename-InlineLogo $section.Footers.Item(2) "image2.png"

# BTec logo, first-page header
Rename-InlineLogo $section.Headers.Item(2) "image1.jpg"
